$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 5) down onto the
# new row 6 so the appended row keeps the same borders/shading as the rest
# of the table.
$ws.Range("A5:E5").Copy()
$ws.Range("A6:E6").PasteSpecial(-4122)

# New "TestCase_F5" record.
$ws.Range("A6").Value = "TestCase_F5"
$ws.Range("B6").Value = "OPQA-877"
$ws.Range("C6").Value = "Verify that user receives a notification when someone he is following  publishes a post"
$ws.Range("D6").Value = "Y"
$ws.Range("E6").Value = "SKIP"

# Match the final selection left behind in the workbook (active cell D5).
$ws.Range("D5").Select()
